$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G4").Value = "2016-09-05 18:51:13"

$wsZhCn.Range("H4").Value = "2016-09-05 18:51:02"
$wsZhCn.Range("K4").Value = "2016-09-05 18:51:31"

$wsDeDe.Range("H4").Value = "2016-09-05 18:51:13"
$wsDeDe.Range("K4").Value = "2016-09-05 18:51:39"
